$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.95578266666667
$ws.Range("H2").Value = 59.867348
$ws.Range("I2").Value = 0.0117373419656925
$ws.Range("J2").Value = 0.0117373419656925
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 13.441269
$ws.Range("N2").Value = 40.323807
$ws.Range("O2").Value = 0.08973082133481231
$ws.Range("P2").Value = 0.08973082133481232
$ws.Range("Q2").Value = 268.231042928204
$ws.Range("R2").Value = 2414.079386353836
$ws.Range("S2").Value = 0.001053201334869148
$ws.Range("T2").Value = 0.001053201334869148
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.95578266666667
$ws.Range("H3").Value = 59.867348
$ws.Range("I3").Value = 0.0117373419656925
$ws.Range("J3").Value = 0.0117373419656925
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 54.711535
$ws.Range("N3").Value = 164.134605
$ws.Range("O3").Value = 0.3652416280068742
$ws.Range("P3").Value = 0.3652416280068742
$ws.Range("Q3").Value = 1091.811501819727
$ws.Range("R3").Value = 9826.30351637754
$ws.Range("S3").Value = 0.004286965888022931
$ws.Range("T3").Value = 0.004286965888022932
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.95578266666667
$ws.Range("H4").Value = 59.867348
$ws.Range("I4").Value = 0.0117373419656925
$ws.Range("J4").Value = 0.0117373419656925
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 63.67711
$ws.Range("N4").Value = 191.03133
$ws.Range("O4").Value = 0.4250937452800914
$ws.Range("P4").Value = 0.4250937452800915
$ws.Range("Q4").Value = 1270.726568001427
$ws.Range("R4").Value = 11436.53911201284
$ws.Range("S4").Value = 0.004989470655829414
$ws.Range("T4").Value = 0.004989470655829414
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.95578266666667
$ws.Range("H5").Value = 59.867348
$ws.Range("I5").Value = 0.0117373419656925
$ws.Range("J5").Value = 0.0117373419656925
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.96553866666667
$ws.Range("N5").Value = 53.896616
$ws.Range("O5").Value = 0.119933805378222
$ws.Range("P5").Value = 0.119933805378222
$ws.Range("Q5").Value = 358.5163851215965
$ws.Range("R5").Value = 3226.647466094368
$ws.Range("S5").Value = 0.001407704086971001
$ws.Range("T5").Value = 0.001407704086971002
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1637.343343333333
$ws.Range("H6").Value = 4912.03003
$ws.Range("I6").Value = 0.9630320723052701
$ws.Range("J6").Value = 0.9630320723052702
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 13.441269
$ws.Range("N6").Value = 40.323807
$ws.Range("O6").Value = 0.08973082133481231
$ws.Range("P6").Value = 0.08973082133481232
$ws.Range("Q6").Value = 22007.97232310269
$ws.Range("R6").Value = 198071.7509079242
$ws.Range("S6").Value = 0.08641365881971824
$ws.Range("T6").Value = 0.08641365881971826
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1637.343343333333
$ws.Range("H7").Value = 4912.03003
$ws.Range("I7").Value = 0.9630320723052701
$ws.Range("J7").Value = 0.9630320723052702
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 54.711535
$ws.Range("N7").Value = 164.134605
$ws.Range("O7").Value = 0.3652416280068742
$ws.Range("P7").Value = 0.3652416280068742
$ws.Range("Q7").Value = 89581.56763579867
$ws.Range("R7").Value = 806234.108722188
$ws.Range("S7").Value = 0.3517394019116106
$ws.Range("T7").Value = 0.3517394019116107
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1637.343343333333
$ws.Range("H8").Value = 4912.03003
$ws.Range("I8").Value = 0.9630320723052701
$ws.Range("J8").Value = 0.9630320723052702
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 63.67711
$ws.Range("N8").Value = 191.03133
$ws.Range("O8").Value = 0.4250937452800914
$ws.Range("P8").Value = 0.4250937452800915
$ws.Range("Q8").Value = 104261.2921812044
$ws.Range("R8").Value = 938351.6296308399
$ws.Range("S8").Value = 0.4093789104410951
$ws.Range("T8").Value = 0.4093789104410951
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1637.343343333333
$ws.Range("H9").Value = 4912.03003
$ws.Range("I9").Value = 0.9630320723052701
$ws.Range("J9").Value = 0.9630320723052702
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.96553866666667
$ws.Range("N9").Value = 53.896616
$ws.Range("O9").Value = 0.119933805378222
$ws.Range("P9").Value = 0.119933805378222
$ws.Range("Q9").Value = 29415.75514526427
$ws.Range("R9").Value = 264741.7963073785
$ws.Range("S9").Value = 0.1155001011328461
$ws.Range("T9").Value = 0.1155001011328461
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 17.50081933333334
$ws.Range("H10").Value = 52.502458
$ws.Range("I10").Value = 0.01029341242216722
$ws.Range("J10").Value = 0.01029341242216722
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.441269
$ws.Range("N10").Value = 40.323807
$ws.Range("O10").Value = 0.08973082133481231
$ws.Range("P10").Value = 0.08973082133481232
$ws.Range("Q10").Value = 235.233220379734
$ws.Range("R10").Value = 2117.098983417606
$ws.Range("S10").Value = 0.0009236363509790242
$ws.Range("T10").Value = 0.0009236363509790245
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 17.50081933333334
$ws.Range("H11").Value = 52.502458
$ws.Range("I11").Value = 0.01029341242216722
$ws.Range("J11").Value = 0.01029341242216722
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 54.711535
$ws.Range("N11").Value = 164.134605
$ws.Range("O11").Value = 0.3652416280068742
$ws.Range("P11").Value = 0.3652416280068742
$ws.Range("Q11").Value = 957.4966894843434
$ws.Range("R11").Value = 8617.47020535909
$ws.Range("S11").Value = 0.003759582710818536
$ws.Range("T11").Value = 0.003759582710818538
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 17.50081933333334
$ws.Range("H12").Value = 52.502458
$ws.Range("I12").Value = 0.01029341242216722
$ws.Range("J12").Value = 0.01029341242216722
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 63.67711
$ws.Range("N12").Value = 191.03133
$ws.Range("O12").Value = 0.4250937452800914
$ws.Range("P12").Value = 0.4250937452800915
$ws.Range("Q12").Value = 1114.401597778794
$ws.Range("R12").Value = 10029.61438000914
$ws.Range("S12").Value = 0.00437566523825168
$ws.Range("T12").Value = 0.004375665238251681
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 17.50081933333334
$ws.Range("H13").Value = 52.502458
$ws.Range("I13").Value = 0.01029341242216722
$ws.Range("J13").Value = 0.01029341242216722
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 17.96553866666667
$ws.Range("N13").Value = 53.896616
$ws.Range("O13").Value = 0.119933805378222
$ws.Range("P13").Value = 0.119933805378222
$ws.Range("Q13").Value = 314.4116464313476
$ws.Range("R13").Value = 2829.704817882128
$ws.Range("S13").Value = 0.001234528122117976
$ws.Range("T13").Value = 0.001234528122117976
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 25.39612333333333
$ws.Range("H14").Value = 76.18836999999999
$ws.Range("I14").Value = 0.01493717330687017
$ws.Range("J14").Value = 0.01493717330687017
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 13.441269
$ws.Range("N14").Value = 40.323807
$ws.Range("O14").Value = 0.08973082133481231
$ws.Range("P14").Value = 0.08973082133481232
$ws.Range("Q14").Value = 341.35612528051
$ws.Range("R14").Value = 3072.20512752459
$ws.Range("S14").Value = 0.001340324829245895
$ws.Range("T14").Value = 0.001340324829245895
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 25.39612333333333
$ws.Range("H15").Value = 76.18836999999999
$ws.Range("I15").Value = 0.01493717330687017
$ws.Range("J15").Value = 0.01493717330687017
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 54.711535
$ws.Range("N15").Value = 164.134605
$ws.Range("O15").Value = 0.3652416280068742
$ws.Range("P15").Value = 0.3652416280068742
$ws.Range("Q15").Value = 1389.460890615983
$ws.Range("R15").Value = 12505.14801554385
$ws.Range("S15").Value = 0.005455677496422084
$ws.Range("T15").Value = 0.005455677496422085
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 25.39612333333333
$ws.Range("H16").Value = 76.18836999999999
$ws.Range("I16").Value = 0.01493717330687017
$ws.Range("J16").Value = 0.01493717330687017
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 63.67711
$ws.Range("N16").Value = 191.03133
$ws.Range("O16").Value = 0.4250937452800914
$ws.Range("P16").Value = 0.4250937452800915
$ws.Range("Q16").Value = 1617.151739070233
$ws.Range("R16").Value = 14554.3656516321
$ws.Range("S16").Value = 0.006349698944915248
$ws.Range("T16").Value = 0.006349698944915248
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 25.39612333333333
$ws.Range("H17").Value = 76.18836999999999
$ws.Range("I17").Value = 0.01493717330687017
$ws.Range("J17").Value = 0.01493717330687017
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 17.96553866666667
$ws.Range("N17").Value = 53.896616
$ws.Range("O17").Value = 0.119933805378222
$ws.Range("P17").Value = 0.119933805378222
$ws.Range("Q17").Value = 456.2550357284355
$ws.Range("R17").Value = 4106.29532155592
$ws.Range("S17").Value = 0.001791472036286939
$ws.Range("T17").Value = 0.00179147203628694

Write-Output "done"
